$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 18:34"

# Update country name order (ranking changes) in column A
$ws.Range("A16").Value = "India"
$ws.Range("A17").Value = "Peru"
$ws.Range("A46").Value = "Republica Dominicana"
$ws.Range("A47").Value = "Sudafrica"
$ws.Range("A55").Value = "Marruecos"
$ws.Range("A56").Value = "Finlandia"
$ws.Range("A95").Value = "Somalia"
$ws.Range("A96").Value = "Consejo Danes para los Refugiados"
$ws.Range("A97").Value = "Kirguistan"
$ws.Range("A98").Value = "Letonia"

# Update numeric statistics per row
# Row 4
$ws.Range("B4").Value = 1327320
$ws.Range("C4").Value = 5535
$ws.Range("D4").Value = 224633
$ws.Range("E4").Value = 1023858
$ws.Range("G4").Value = 214
$ws.Range("H4").Value = 78829
# Row 6
$ws.Range("B6").Value = 218268
$ws.Range("C6").Value = 1083
$ws.Range("D6").Value = 103031
$ws.Range("E6").Value = 84842
$ws.Range("F6").Value = 1034
$ws.Range("G6").Value = 194
$ws.Range("H6").Value = 30395
# Row 11
$ws.Range("B11").Value = 147261
$ws.Range("C11").Value = 1369
$ws.Range("E11").Value = 77920
$ws.Range("G11").Value = 52
$ws.Range("H11").Value = 10044
# Row 15
$ws.Range("B15").Value = 66780
$ws.Range("C15").Value = 346
$ws.Range("E15").Value = 31746
$ws.Range("G15").Value = 59
$ws.Range("H15").Value = 4628
# Row 16
$ws.Range("B16").Value = 62521
$ws.Range("C16").Value = 2826
$ws.Range("D16").Value = 18672
$ws.Range("E16").Value = 41760
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 104
$ws.Range("H16").Value = 2089
# Row 17
$ws.Range("B17").Value = 61847
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 19012
$ws.Range("E17").Value = 41121
$ws.Range("F17").Value = 730
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1714
# Row 46
$ws.Range("B46").Value = 9882
$ws.Range("C46").Value = 506
$ws.Range("D46").Value = 2584
$ws.Range("E46").Value = 6913
$ws.Range("F46").Value = 134
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 385
# Row 47
$ws.Range("B47").Value = 9420
$ws.Range("C47").Value = 525
$ws.Range("D47").Value = 3983
$ws.Range("E47").Value = 5251
$ws.Range("F47").Value = 77
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 186
# Row 55
$ws.Range("B55").Value = 5910
$ws.Range("C55").Value = 199
$ws.Range("D55").Value = 2461
$ws.Range("E55").Value = 3263
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 186
# Row 56
$ws.Range("B56").Value = 5880
$ws.Range("C56").Value = 142
$ws.Range("D56").Value = 4000
$ws.Range("E56").Value = 1615
$ws.Range("F56").Value = 45
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 265
# Row 57
$ws.Range("D57").Value = 1728
$ws.Range("E57").Value = 3590
# Row 71
$ws.Range("B71").Value = 2679
$ws.Range("C71").Value = 76
$ws.Range("D71").Value = 1702
$ws.Range("E71").Value = 870
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 107
# Row 95
$ws.Range("B95").Value = 997
$ws.Range("C95").Value = 69
$ws.Range("D95").Value = 110
$ws.Range("E95").Value = 839
$ws.Range("F95").Value = 2
$ws.Range("G95").Value = 4
$ws.Range("H95").Value = 48
# Row 96
$ws.Range("B96").Value = 937
$ws.Range("C96").Value = 74
$ws.Range("D96").Value = 130
$ws.Range("E96").Value = 768
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 39
# Row 97
$ws.Range("B97").Value = 931
$ws.Range("C97").Value = 25
$ws.Range("D97").Value = 658
$ws.Range("E97").Value = 261
$ws.Range("F97").Value = 13
$ws.Range("H97").Value = 12
# Row 98
$ws.Range("B98").Value = 930
$ws.Range("C98").Value = 2
$ws.Range("D98").Value = 464
$ws.Range("E98").Value = 448
$ws.Range("H98").Value = 18
# Row 109
$ws.Range("B109").Value = 754
$ws.Range("C109").Value = 2
$ws.Range("D109").Value = 545
$ws.Range("E109").Value = 161
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 48
